# Questions.docx edit:
#  - Remove the "Customer retention" and "Profit" + "ability" bullet
#    paragraphs entirely.
#  - Remove the (now redundant) empty ListParagraph that carried the
#    _GoBack bookmark, folding that bookmark onto the end of the
#    "Year wise sales" bullet instead.
#  - On the final (empty) paragraph, swap the ListParagraph style for an
#    explicit left indent of 360 twips.

$d = $word.ActiveDocument

# --- Step 1: delete the "Customer retention" and "Profit"+"ability" paragraphs ---
$custParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Customer retention") {
        $custParaIndex = $i
        break
    }
}
$profitParaIndex = $custParaIndex + 1

$delRange = $d.Range($d.Paragraphs($custParaIndex).Range.Start, $d.Paragraphs($profitParaIndex).Range.End)
$delRange.Delete()

# --- Step 2: find the "Year wise sales" paragraph and the empty bookmark paragraph right after it ---
$yearParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Year wise sales") {
        $yearParaIndex = $i
        break
    }
}

$yearPara = $d.Paragraphs($yearParaIndex)

# Re-write the "Year wise sales" paragraph so it keeps its own (ListParagraph +
# numbered list) formatting but also carries the _GoBack bookmark at its end.
$yearXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Year wise sales</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$yearPara.Range.InsertXML($yearXml)

# The paragraph right after "Year wise sales" is now the empty ListParagraph
# that used to hold the bookmark; it is redundant, delete it outright.
$oldBookmarkPara = $d.Paragraphs($yearParaIndex + 1)
$oldBookmarkPara.Range.Delete()

# --- Step 3: fix up the final (empty) paragraph's formatting ---
# Replace "pStyle=ListParagraph" with an explicit "ind left=360" on the very
# last paragraph. Because that paragraph's Range is just its own end-of-
# paragraph mark, InsertXML on it alone would insert a sibling rather than
# replace it in place, so instead we replace the whole span from the end of
# the (now bookmarked) "Year wise sales" paragraph through the end of the
# document, keeping "Year wise sales" untouched and only emitting the new
# trailing paragraph.
$yearPara = $d.Paragraphs($yearParaIndex)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tailRange = $d.Range($yearPara.Range.End, $lastPara.Range.End)

$tailXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$tailRange.InsertXML($tailXml)

# The InsertXML above inserts the new trailing paragraph before the old
# (still ListParagraph-styled) final paragraph rather than replacing it in
# place, so delete the now-duplicate original final paragraph by removing
# its end-of-paragraph mark together with the rest of the tail.
$newTailPara = $d.Paragraphs($yearParaIndex + 1)
$dupPara = $d.Paragraphs($d.Paragraphs.Count)
if ($dupPara.Range.Start -ne $newTailPara.Range.Start) {
    $cleanupRange = $d.Range($newTailPara.Range.End - 1, $dupPara.Range.End)
    $cleanupRange.Delete()
}

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i : $($d.Paragraphs($i).Range.Text)"
}
